$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column M mirrors column L's formatting (the 2021 column), extended for 2022.
$ws.Range("L3").Copy()
$ws.Range("M3").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("L4").Copy()
$ws.Range("M4").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("M4").Value = 2022

$ws.Range("L6").Copy()
$ws.Range("M6").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("M6").Value = 18

$ws.Range("L7").Copy()
$ws.Range("M7").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("M7").Value = 6.2

$ws.Range("L8").Copy()
$ws.Range("M8").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("M8").Value = "-"

$excel.CutCopyMode = 0

# Move the active selection from N5 to N4, as in the edited workbook.
$ws.Range("N4").Select()
